# Apply cryptocurrency price/volume updates as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "64.219.98" },
    @{ Cell = "E2"; Value = "  -1.78%  " },
    @{ Cell = "D3"; Value = "3.118.72" },
    @{ Cell = "E3"; Value = "  -2.59%  " },
    @{ Cell = "E4"; Value = "  -0.02%  " },
    @{ Cell = "D5"; Value = "594.57" },
    @{ Cell = "E5"; Value = "  -0.60%  " },
    @{ Cell = "D6"; Value = "158.03" },
    @{ Cell = "E6"; Value = "  +2.84%  " },
    @{ Cell = "E7"; Value = "  +0.00%  " },
    @{ Cell = "D8"; Value = "0.542" },
    @{ Cell = "E8"; Value = "  -0.03%  " },
    @{ Cell = "D9"; Value = "3.117.05" },
    @{ Cell = "E9"; Value = "  -2.59%  " },
    @{ Cell = "E10"; Value = "  -5.52%  " },
    @{ Cell = "D11"; Value = "5.94" },
    @{ Cell = "E11"; Value = "  -2.89%  " },
    @{ Cell = "D12"; Value = "0.453" },
    @{ Cell = "E12"; Value = "  -4.05%  " },
    @{ Cell = "D13"; Value = "37.25" },
    @{ Cell = "E13"; Value = "  -5.32%  " },
    @{ Cell = "D14"; Value = "0.0000240" },
    @{ Cell = "E14"; Value = "  -5.86%  " },
    @{ Cell = "D15"; Value = "3.638.43" },
    @{ Cell = "E15"; Value = "  -2.40%  " },
    @{ Cell = "E16"; Value = "  -1.40%  " },
    @{ Cell = "D17"; Value = "7.27" },
    @{ Cell = "E17"; Value = "  -2.13%  " },
    @{ Cell = "D18"; Value = "64.175.64" },
    @{ Cell = "E18"; Value = "  -1.39%  " },
    @{ Cell = "D19"; Value = "3.122.11" },
    @{ Cell = "E19"; Value = "  -2.45%  " },
    @{ Cell = "D20"; Value = "478.40" },
    @{ Cell = "E20"; Value = "  -1.15%  " },
    @{ Cell = "D21"; Value = "14.56" },
    @{ Cell = "E21"; Value = "  -3.64%  " },
    @{ Cell = "D22"; Value = "0.717" },
    @{ Cell = "E22"; Value = "  -7.55%  " },
    @{ Cell = "D23"; Value = "7.59" },
    @{ Cell = "E23"; Value = "  -4.38%  " },
    @{ Cell = "E24"; Value = "  +1.22%  " },
    @{ Cell = "D25"; Value = "12.98" },
    @{ Cell = "E25"; Value = "  -6.89%  " },
    @{ Cell = "D26"; Value = "81.45" },
    @{ Cell = "E26"; Value = "  -2.67%  " },
    @{ Cell = "D27"; Value = "10.51" },
    @{ Cell = "E27"; Value = "  +6.85%  " },
    @{ Cell = "E28"; Value = "  -0.29%  " },
    @{ Cell = "D29"; Value = "7.63" },
    @{ Cell = "E29"; Value = "  +1.92%  " },
    @{ Cell = "D30"; Value = "2.70" },
    @{ Cell = "E30"; Value = "  -2.84%  " },
    @{ Cell = "B31"; Value = "FirstDigitalUSD" },
    @{ Cell = "C31"; Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd" },
    @{ Cell = "D31"; Value = "1.00" },
    @{ Cell = "E31"; Value = "  -0.12%  " },
    @{ Cell = "B32"; Value = "ImmutableX" },
    @{ Cell = "C32"; Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" },
    @{ Cell = "D32"; Value = "2.21" },
    @{ Cell = "E32"; Value = "  -3.06%  " },
    @{ Cell = "D33"; Value = "0.113" },
    @{ Cell = "E33"; Value = "  -6.23%  " },
    @{ Cell = "D34"; Value = "27.36" },
    @{ Cell = "E34"; Value = "  -4.39%  " },
    @{ Cell = "D35"; Value = "0.0₃0848" },
    @{ Cell = "E35"; Value = "  -5.62%  " },
    @{ Cell = "E36"; Value = "  -2.76%  " },
    @{ Cell = "B37"; Value = "Filecoin" },
    @{ Cell = "C37"; Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil" },
    @{ Cell = "D37"; Value = "6.05" },
    @{ Cell = "E37"; Value = "  -4.86%  " },
    @{ Cell = "B38"; Value = "dogwifhat" },
    @{ Cell = "C38"; Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif" },
    @{ Cell = "D38"; Value = "3.32" },
    @{ Cell = "E38"; Value = "  -7.59%  " },
    @{ Cell = "E39"; Value = "  -5.52%  " },
    @{ Cell = "D40"; Value = "51.12" },
    @{ Cell = "E40"; Value = "  -1.13%  " },
    @{ Cell = "D41"; Value = "9.17" },
    @{ Cell = "E41"; Value = "  -3.13%  " },
    @{ Cell = "D42"; Value = "448.64" },
    @{ Cell = "E42"; Value = "  -5.41%  " },
    @{ Cell = "D43"; Value = "0.292" },
    @{ Cell = "E43"; Value = "  -3.42%  " },
    @{ Cell = "D44"; Value = "0.0366" },
    @{ Cell = "E44"; Value = "  -4.73%  " },
    @{ Cell = "D45"; Value = "0.112" },
    @{ Cell = "E45"; Value = "  +0.15%  " },
    @{ Cell = "D46"; Value = "40.30" },
    @{ Cell = "E46"; Value = "  +4.13%  " },
    @{ Cell = "D47"; Value = "2.832.81" },
    @{ Cell = "E47"; Value = "  -4.33%  " },
    @{ Cell = "D48"; Value = "130.84" },
    @{ Cell = "E48"; Value = "  -0.62%  " },
    @{ Cell = "D49"; Value = "25.81" },
    @{ Cell = "E49"; Value = "  +0.64%  " },
    @{ Cell = "E50"; Value = "  +0.01%  " },
    @{ Cell = "D51"; Value = "2.26" },
    @{ Cell = "E51"; Value = "  -3.30%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (e.g. "594.57") are not
    # coerced into floating point numbers, preserving the original inline/shared
    # string representation and avoiding precision drift.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}

